$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected Price/Volume cells to stay text (avoid Excel
# auto-parsing values like "25.65" or "1.000" as numbers), matching
# the original inlineStr cell content, then restore default styling.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.876.43'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.888.00'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '0.7706'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '242.62'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("D9").Value = '25.65'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").Value = '0.07164'
$ws.Range("E10").Value = '  -5.88%  '
$ws.Range("D11").Value = '0.08592'
$ws.Range("E11").Value = '  +5.86%  '
$ws.Range("D12").Value = '0.7645'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").Value = '1.899.31'
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '5.364'
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("D15").Value = '93.65'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").Value = '6.154'
$ws.Range("E16").Value = '  -1.43%  '
$ws.Range("D17").Value = '29.925.86'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = '13.77'
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '0.000007804'
$ws.Range("D21").Value = '2.174.68'
$ws.Range("E21").Value = '  +5.66%  '
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '8.043'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").Value = '0.1638'
$ws.Range("E25").Value = '  +4.40%  '
$ws.Range("D26").Value = '9.380'
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("D27").Value = '162.60'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").Value = '4.511'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '4.096'
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").Value = '0.05466'
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").Value = '1.241'
$ws.Range("E35").Value = '  -1.79%  '
$ws.Range("D36").Value = '0.7443'
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("D37").Value = '1.002'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '2.699'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("D39").Value = '0.01952'
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("D40").Value = '2.781'
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("D41").Value = '0.4469'
$ws.Range("D42").Value = '1.112.92'
$ws.Range("E42").Value = '  -4.12%  '
$ws.Range("D43").Value = '6.089'
$ws.Range("E43").Value = '  +2.32%  '
$ws.Range("D44").Value = '73.03'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = '0.8518'
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '102.50'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").Value = '1.864'
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  +1.23%  '
$ws.Range("E50").Value = '  -3.54%  '
$ws.Range("D51").Value = '2.067.71'
$ws.Range("E51").Value = '  +2.05%  '

# Restore default (unstyled) cell formatting so the saved cells match
# the workbook's original plain text-cell appearance.
$dataRange.Style = "Normal"
